# "final project idea added"
#
# Append a new bulleted "Gift Search" idea paragraph (as two runs of text,
# mirroring the source), followed by a new empty bulleted paragraph, to the
# existing ListParagraph / numId=1 bullet list at the end of the document.

$d = $word.ActiveDocument

# The last existing paragraph ("The Invisible Stalker: ...") already carries
# the ListParagraph style + numPr (ilvl 0, numId 1) that the new paragraphs
# should inherit, so create the new paragraphs right after it - Word
# automatically carries the paragraph-level formatting forward.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# --- New paragraph: the "Gift Search" idea ---------------------------------
$giftPara = $d.Paragraphs.Last
$firstSentence = "Gift Search: a program that helps find the right website for a particular gift. The person will be asked various questions about the type of gift they are looking for such as toy, electronic, clothing, shoes, or other. Then possibly age range followed by price range"
$secondSentence = " which will give a close estimate of which area and maybe even which store they are needing."
$giftPara.Range.InsertAfter($firstSentence + $secondSentence)

# --- New trailing empty bulleted paragraph ----------------------------------
$giftPara.Range.InsertParagraphAfter()

# Keep the "Gift Search" text as two distinct runs (matching how it was
# authored as two separate sentences) instead of letting them collapse into
# one merged run, by re-touching the character formatting of the second
# sentence only - a momentary flip (on, then back off) is enough for Word to
# record it as its own run without changing the visible formatting.
$giftText = $giftPara.Range.Text
$splitIndex = $giftText.IndexOf($secondSentence)
if ($splitIndex -ge 0) {
    $splitStart = $giftPara.Range.Start + $splitIndex
    $splitEnd = $giftPara.Range.End - 1
    $secondRun = $d.Range($splitStart, $splitEnd)
    $secondRun.Font.Bold = 1
    $secondRun.Font.Bold = 0
}

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
Write-Output ("GiftParagraph=" + $giftPara.Range.Text)
